{"js": "// Update the division-problem worksheet table: each populated cell's\n// expression text is replaced with a new one, per the commit's data.\n// Cells are addressed by (row, column) position in the table rather than by\n// searching for the old text, because several of the new values coincide\n// with *other* cells' original values (e.g. \"23\u00f73=\" and \"83\u00f75=\" are both a\n// pre-existing cell's text AND another cell's replacement text), which would\n// make a global text search-and-replace ambiguous/unsafe once earlier\n// replacements have been applied.\nconst replacements = [\n  { row: 0, cells: [\"23\u00f73=\", \"76\u00f76=\", \"29\u00f79=\", \"64\u00f75=\", \"17\u00f76=\"] },\n  { row: 4, cells: [\"92\u00f76=\", \"22\u00f72=\", \"42\u00f74=\", \"49\u00f72=\", \"27\u00f77=\"] },\n  { row: 8, cells: [\"58\u00f78=\", \"24\u00f75=\", \"14\u00f77=\", \"83\u00f72=\", \"53\u00f77=\"] },\n  { row: 12, cells: [\"16\u00f72=\", \"83\u00f75=\", \"76\u00f77=\", \"93\u00f77=\", \"92\u00f74=\"] },\n  { row: 16, cells: [\"98\u00f76=\", \"37\u00f74=\", \"48\u00f79=\", \"10\u00f72=\", \"50\u00f76=\"] },\n];\n\nconst table = context.document.body.tables.getFirst();\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\n// Gather the paragraph we'll rewrite for every target cell first.\nconst targets = [];\nfor (const { row, cells: newTexts } of replacements) {\n  const rowCells = rows.items[row].cells;\n  rowCells.load(\"items\");\n  targets.push({ rowCells, newTexts });\n}\nawait context.sync();\n\nconst paragraphsToWrite = [];\nfor (const { rowCells, newTexts } of targets) {\n  for (let col = 0; col < newTexts.length; col++) {\n    const cell = rowCells.items[col];\n    const paragraph = cell.body.paragraphs.getFirst();\n    paragraphsToWrite.push({ paragraph, newText: newTexts[col] });\n  }\n}\n\nfor (const { paragraph, newText } of paragraphsToWrite) {\n  // Replacing at the paragraph level (rather than the cell body) keeps the\n  // existing run/paragraph formatting (font, size, alignment) intact.\n  paragraph.insertText(newText, Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Update the division-problem worksheet table: each populated cell's\n# expression text is replaced with a new one, per the commit's data.\n# Cells are addressed by (row, column) position in the table rather than by\n# searching for the old text, because several of the new values coincide\n# with *other* cells' original values (e.g. \"23\u00f73=\" and \"83\u00f75=\" are both a\n# pre-existing cell's text AND another cell's replacement text), which would\n# make a global Find/Replace ambiguous/unsafe once earlier replacements have\n# already been applied.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$rowReplacements = @{\n    1  = @(\"23\u00f73=\", \"76\u00f76=\", \"29\u00f79=\", \"64\u00f75=\", \"17\u00f76=\")\n    5  = @(\"92\u00f76=\", \"22\u00f72=\", \"42\u00f74=\", \"49\u00f72=\", \"27\u00f77=\")\n    9  = @(\"58\u00f78=\", \"24\u00f75=\", \"14\u00f77=\", \"83\u00f72=\", \"53\u00f77=\")\n    13 = @(\"16\u00f72=\", \"83\u00f75=\", \"76\u00f77=\", \"93\u00f77=\", \"92\u00f74=\")\n    17 = @(\"98\u00f76=\", \"37\u00f74=\", \"48\u00f79=\", \"10\u00f72=\", \"50\u00f76=\")\n}\n\nforeach ($row in $rowReplacements.Keys) {\n    $newTexts = $rowReplacements[$row]\n    for ($col = 1; $col -le $newTexts.Length; $col++) {\n        $cell = $t.Cell($row, $col)\n        # Assigning Range.Text replaces just the cell's text content while\n        # keeping the existing run/paragraph formatting (font, size,\n        # alignment) intact (Word auto-preserves the trailing cell mark).\n        $cell.Range.Text = $newTexts[$col - 1]\n    }\n}\n"}
